$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in A1 (new data pull at 20:21)
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 20:21"

# Row data: row number, country name (only set if it changed position/name), then B..H values
$rows = @(
    @(4, "Estados Unidos", 7019199, 14431, 4272199, 2542754, 0, 128, 204246),
    @(5, "India", 5557517, 71905, 4492134, 976440, 0, 1034, 88943),
    @(11, "España", 671468, 2957, 0, 0, 0, 56, 30663),
    @(12, "Sudafrica", 661211, 0, 590071, 55187, 0, 0, 15953),
    @(14, "Francia", 458061, 5298, 91574, 335149, 0, 53, 31338),
    @(25, "Alemania", 274997, 1520, 244000, 21520, 0, 7, 9477),
    @(27, "Israel", 190037, 2135, 136502, 52263, 0, 16, 1272),
    @(29, "Canada", 144664, 1015, 125204, 10239, 0, 4, 9221),
    @(37, "Marruecos", 103119, 1376, 84158, 17106, 0, 25, 1855),
    @(38, "Belgica", 102295, 1547, 18965, 73382, 0, 4, 9948),
    @(39, "Egipto", 102015, 0, 89532, 6713, 0, 0, 5770),
    @(70, "Kenia", 37079, 98, 23949, 12480, 0, 2, 650),
    @(73, "Irlanda", 33121, 188, 23364, 7965, 0, 0, 1792),
    @(74, "Serbia", 32938, 30, 31536, 659, 0, 2, 743),
    @(100, "Maldivas", 9770, 46, 8390, 1346, 0, 1, 34),
    @(142, "Sri Lanka", 3299, 12, 3100, 186, 0, 0, 13),
    @(189, "Islas Caimanes", 209, 1, 204, 4, 0, 0, 1),
    @(194, "Seychelles", 143, 2, 136, 7, 0, 0, 0),
    @(204, "Timor Oriental", 27, 0, 26, 1, 0, 0, 0),
    @(205, "Santa Lucia", 27, 0, 26, 1, 0, 0, 0),
    @(211, "San Cristobal y Nieves", 19, 2, 17, 2, 0, 0, 0),
    @(212, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $countryName = $entry[1]
    $ws.Cells.Item($r, 1).Value = $countryName
    for ($i = 0; $i -lt 7; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $entry[2 + $i]
    }
}